$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.602.49'
$ws.Range("E2").Value = '  +1.22%  '
# Row 3
$ws.Range("D3").Value = '3.450.40'
$ws.Range("E3").Value = '  +2.37%  '
# Row 4
$ws.Range("E4").Value = '  -0.13%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.47'
$ws.Range("E5").Value = '  +1.74%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.36'
$ws.Range("E6").Value = '  +5.92%  '
# Row 7
$ws.Range("D7").Value = '3.451.96'
$ws.Range("E7").Value = '  +2.49%  '
# Row 8
$ws.Range("E8").Value = '  +0.01%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.477'
$ws.Range("E9").Value = '  +2.13%  '
# Row 10
$ws.Range("E10").Value = '  +0.31%  '
# Row 11
$ws.Range("E11").Value = '  +3.61%  '
# Row 12
$ws.Range("E12").Value = '  +2.70%  '
# Row 13
$ws.Range("D13").Value = '4.039.82'
$ws.Range("E13").Value = '  +2.40%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.99'
$ws.Range("E14").Value = '  +9.48%  '
# Row 15
$ws.Range("E15").Value = '  -0.85%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  +2.13%  '
# Row 17
$ws.Range("D17").Value = '3.446.51'
$ws.Range("E17").Value = '  +2.16%  '
# Row 18
$ws.Range("D18").Value = '61.693.70'
$ws.Range("E18").Value = '  +1.03%  '
# Row 19
$ws.Range("E19").Value = '  +9.41%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.28'
$ws.Range("E20").Value = '  +3.97%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.53'
$ws.Range("E21").Value = '  +2.43%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.25'
$ws.Range("E22").Value = '  +4.40%  '
# Row 23
$ws.Range("E23").Value = '  +3.57%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.43'
$ws.Range("E24").Value = '  +3.61%  '
# Row 25
$ws.Range("E25").Value = '  +0.28%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.21%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").Value = '  +1.26%  '
# Row 28
$ws.Range("D28").Value = '3.588.42'
$ws.Range("E28").Value = '  +2.10%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.180'
$ws.Range("E29").Value = '  +1.56%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.60'
$ws.Range("E30").Value = '  +4.09%  '
# Row 31
$ws.Range("E31").Value = '  -0.05%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  +1.81%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("E33").Value = '  +3.02%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").Value = '  -10.99%  '
# Row 35
$ws.Range("E35").Value = '  -0.03%  '
# Row 36
$ws.Range("E36").Value = '  +3.51%  '
# Row 37
$ws.Range("D37").Value = '3.477.92'
$ws.Range("E37").Value = '  +2.49%  '
# Row 38
$ws.Range("E38").Value = '  +3.56%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.13'
$ws.Range("E39").Value = '  +0.50%  '
# Row 40
$ws.Range("E40").Value = '  +1.18%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '166.73'
$ws.Range("E41").Value = '  +1.19%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '28.21'
$ws.Range("E42").Value = '  +14.62%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0783'
$ws.Range("E43").Value = '  +3.64%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.804'
$ws.Range("E44").Value = '  +4.05%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.12%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.36'
$ws.Range("E46").Value = '  +1.81%  '
# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.73'
$ws.Range("E47").Value = '  +3.41%  '
# Row 48
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.48'
$ws.Range("E48").Value = '  +4.16%  '
# Row 49
$ws.Range("D49").Value = '2.587.71'
$ws.Range("E49").Value = '  +1.92%  '
# Row 50
$ws.Range("E50").Value = '  -1.21%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.91'
$ws.Range("E51").Value = '  +2.46%  '
